# Auto-generated files on 2025-12-16
# Update the HotStock_Top20 ranking table (rows 2-21, columns A:C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = "航天电子"; "B2" = "永辉超市"; "C2" = "航天发展";
    "A3" = "航天发展"; "B3" = "航天发展"; "C3" = "航天电子";
    "A4" = "通宇通讯"; "B4" = "航天电子"; "C4" = "永辉超市";
    "A5" = "永辉超市"; "B5" = "平潭发展"; "C5" = "平潭发展";
    "A6" = "北汽蓝谷"; "B6" = "北汽蓝谷"; "C6" = "东百集团";
    "A7" = "航天机电"; "B7" = "雪人集团"; "C7" = "雪人集团";
    "A8" = "平潭发展"; "B8" = "雷科防务"; "C8" = "太阳电缆";
    "A9" = "华菱线缆"; "B9" = "东百集团"; "C9" = "通宇通讯";
    "A10" = "雪人集团"; "B10" = "通宇通讯"; "C10" = "龙洲股份";
    "A11" = "东百集团"; "B11" = "航天动力"; "C11" = "中超控股";
    "A12" = "太阳电缆"; "B12" = "太阳电缆"; "C12" = "航天信息";
    "A13" = "航天信息"; "B13" = "航天信息"; "C13" = "航天机电";
    "A14" = "雷科防务"; "B14" = "百大集团"; "C14" = "顺灏股份";
    "A15" = "航天动力"; "B15" = "华菱线缆"; "C15" = "海马汽车";
    "A16" = "长安汽车"; "B16" = "航天机电"; "C16" = "华菱线缆";
    "A17" = "百大集团"; "B17" = "长安汽车"; "C17" = "恒宝股份";
    "A18" = "中超控股"; "B18" = "恒宝股份"; "C18" = "航天动力";
    "A19" = "龙洲股份"; "B19" = "摩尔线程-U"; "C19" = "百大集团";
    "A20" = "恒宝股份"; "B20" = "航天科技"; "C20" = "安妮股份";
    "A21" = "顺灏股份"; "B21" = "中超控股"; "C21" = "雷科防务";
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
